$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rng = $ws.Range("A2:D2")
$rng.NumberFormat = "@"
$rng.Value = "10"
